# Apply the content edits described by the commit:
#   "added sound test + bug correction on dual N-back (instructions not
#    shown in the last block) + questionnaires presession 1-3"
#
# Concretely, the OOXML diff shows two deliberate text edits inside the
# "Arrows will move on ..." instruction textbox (shape id 22, present on
# both the English and French instruction slides):
#   EN: "Arrows will move on indicators:"            -> "Arrows will move on gauges:"
#   FR: "Des flèches vont bouger sur des indicateurs:" -> "Des flèches vont bouger sur des jauges:"
#
# Everything else in the diff (ppt/revisionInfo.xml removal, the
# pc:docChgLst block appended to ppt/changesInfos/changesInfo1.xml, and the
# cached "datetimeFigureOut" field text going from 1/31/2025 to 2/17/2025
# across the slide master / slide layouts) is bookkeeping that PowerPoint
# itself stamps when the deck is re-saved on a later day - it is not
# reachable through the PowerPoint object model. We still refresh the
# cached date-field text below, on a best effort basis, since it is the
# only externally visible trace of that resave that the COM object model
# can actually reach.

$p = $ppt.ActivePresentation

function Set-Substring {
    param($TextRange, [string]$Old, [string]$New)

    $full = $TextRange.Text
    $idx = $full.IndexOf($Old)
    if ($idx -ge 0) {
        $sub = $TextRange.Characters($idx + 1, $Old.Length)
        $sub.Text = $New
    }
}

function Update-Indicators {
    param($Shapes)

    for ($i = 1; $i -le $Shapes.Count; $i++) {
        $sh = $Shapes.Item($i)
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
            $tr = $sh.TextFrame.TextRange
            $full = $tr.Text
            if ($full.IndexOf("Arrows will move on indicators:") -ge 0) {
                Set-Substring $tr "Arrows will move on indicators:" "Arrows will move on gauges:"
            }
            elseif ($full.IndexOf("Des flèches vont bouger sur des indicateurs:") -ge 0) {
                Set-Substring $tr "Des flèches vont bouger sur des indicateurs:" "Des flèches vont bouger sur des jauges:"
            }
        }
    }
}

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    Update-Indicators $s.Shapes
}

# --- Best-effort refresh of the cached "datetimeFigureOut" field text ---
$oldDate = "1/31/2025"
$newDate = "2/17/2025"

function Update-DateShape {
    param($Shapes)

    for ($i = 1; $i -le $Shapes.Count; $i++) {
        $sh = $Shapes.Item($i)
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateShape $master.Shapes
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DateShape $layouts.Item($li).Shapes
}
